$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45186 -> 2023-09-17)
# that was bumped by 2 days to 45188 (-> 2023-09-19) for every data row.
$ws.Range("C2:C482").Value = 45188
